$wb = $excel.ActiveWorkbook

# --- Contracts sheet: add ProjectManager, ProjectType, LastUpdatedDate columns ---
$wsContracts = $wb.Worksheets.Item("Contracts")
$wsContracts.Range("I1").Value = "ProjectManager"
$wsContracts.Range("J1").Value = "ProjectType"
$wsContracts.Range("K1").Value = "LastUpdatedDate"

# --- Activities sheet: add ActivityType, DueDate columns ---
$wsActivities = $wb.Worksheets.Item("Activities")
$wsActivities.Range("H1").Value = "ActivityType"
$wsActivities.Range("I1").Value = "DueDate"

# --- Documents sheet: add DocumentDate, Tags columns ---
$wsDocuments = $wb.Worksheets.Item("Documents")
$wsDocuments.Range("K1").Value = "DocumentDate"
$wsDocuments.Range("L1").Value = "Tags"

# --- Add new "Media" sheet at the end of the workbook ---
$wsMedia = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMedia.Name = "Media"
$wsMedia.Range("A1").Value = "ContractID"
$wsMedia.Range("B1").Value = "ImageURL"
$wsMedia.Range("C1").Value = "DateTaken"
$wsMedia.Range("D1").Value = "Description"
$wsMedia.Range("E1").Value = "ActivityID"

# Restore "Documents" as the active/selected sheet (it was active before the edit)
$wsDocuments.Activate()
